$d = $word.ActiveDocument

# Update the header date line (unique text in the document, safe to Find/Replace)
$d.Content.Find.Execute("2025-03-29 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-03-30 Sunday", 2)

# Update the practice-problem table. Addressing cells directly by (row, column)
# avoids any ambiguity from duplicate/overlapping text values that Find/Replace
# could mismatch (e.g. several cells end up sharing the same "333÷9=" text).
$t = $d.Tables.Item(1)

$rowsData = @{
    1  = @("615÷7=", "908÷8=", "214÷8=", "732÷9=", "874÷7=")
    5  = @("746÷8=", "333÷9=", "694÷9=", "862÷2=", "194÷8=")
    9  = @("549÷3=", "995÷5=", "197÷7=", "431÷2=", "841÷8=")
    13 = @("364÷8=", "334÷5=", "798÷5=", "559÷3=", "814÷3=")
    17 = @("881÷3=", "473÷8=", "833÷9=", "570÷4=", "932÷3=")
}

foreach ($rowIndex in $rowsData.Keys) {
    $values = $rowsData[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]
    }
}
